$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 1102
    $ws.Range("F6").Value = 943
}
